$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 197.5433703333333
$ws.Range("H2").Value = 592.6301109999999
$ws.Range("I2").Value = 0.3388703761585983
$ws.Range("J2").Value = 0.3388703761585982
$ws.Range("M2").Value = 11.81073566666667
$ws.Range("N2").Value = 35.432207
$ws.Range("O2").Value = 0.3076347070004043
$ws.Range("P2").Value = 0.3076347070004043
$ws.Range("Q2").Value = 2333.132529709442
$ws.Range("R2").Value = 20998.19276738498
$ws.Range("S2").Value = 0.1042482888806672
$ws.Range("T2").Value = 0.1042482888806672
$ws.Range("G3").Value = 197.5433703333333
$ws.Range("H3").Value = 592.6301109999999
$ws.Range("I3").Value = 0.3388703761585983
$ws.Range("J3").Value = 0.3388703761585982
$ws.Range("O3").Value = 0.3244871420261927
$ws.Range("P3").Value = 0.3244871420261927
$ws.Range("Q3").Value = 2460.943090315108
$ws.Range("R3").Value = 22148.48781283597
$ws.Range("S3").Value = 0.1099590798770444
$ws.Range("T3").Value = 0.1099590798770444
$ws.Range("G4").Value = 197.5433703333333
$ws.Range("H4").Value = 592.6301109999999
$ws.Range("I4").Value = 0.3388703761585983
$ws.Range("J4").Value = 0.3388703761585982
$ws.Range("M4").Value = 3.197710666666667
$ws.Range("N4").Value = 9.593132000000001
$ws.Range("O4").Value = 0.08329089836363292
$ws.Range("P4").Value = 0.0832908983636329
$ws.Range("Q4").Value = 631.6865424441836
$ws.Range("R4").Value = 5685.178881997652
$ws.Range("S4").Value = 0.02822481805907186
$ws.Range("T4").Value = 0.02822481805907185
$ws.Range("G5").Value = 197.5433703333333
$ws.Range("H5").Value = 592.6301109999999
$ws.Range("I5").Value = 0.3388703761585983
$ws.Range("J5").Value = 0.3388703761585982
$ws.Range("M5").Value = 8.081220666666667
$ws.Range("N5").Value = 24.243662
$ws.Range("O5").Value = 0.2104918797744333
$ws.Range("P5").Value = 0.2104918797744333
$ws.Range("Q5").Value = 1596.39156690072
$ws.Range("R5").Value = 14367.52410210648
$ws.Range("S5").Value = 0.07132946247749267
$ws.Range("T5").Value = 0.07132946247749265
$ws.Range("G6").Value = 197.5433703333333
$ws.Range("H6").Value = 592.6301109999999
$ws.Range("I6").Value = 0.3388703761585983
$ws.Range("J6").Value = 0.3388703761585982
$ws.Range("M6").Value = 2.844675333333333
$ws.Range("N6").Value = 8.534026000000001
$ws.Range("O6").Value = 0.07409537283533685
$ws.Range("P6").Value = 0.07409537283533686
$ws.Range("Q6").Value = 561.9467528507651
$ws.Range("R6").Value = 5057.520775656886
$ws.Range("S6").Value = 0.02510872686432218
$ws.Range("T6").Value = 0.02510872686432218
$ws.Range("I7").Value = 0.1369374790620155
$ws.Range("J7").Value = 0.1369374790620154
$ws.Range("M7").Value = 11.81073566666667
$ws.Range("N7").Value = 35.432207
$ws.Range("O7").Value = 0.3076347070004043
$ws.Range("P7").Value = 0.3076347070004043
$ws.Range("Q7").Value = 942.8185802422117
$ws.Range("R7").Value = 8485.367222179904
$ws.Range("S7").Value = 0.04212672124861713
$ws.Range("T7").Value = 0.04212672124861712
$ws.Range("I8").Value = 0.1369374790620155
$ws.Range("J8").Value = 0.1369374790620154
$ws.Range("O8").Value = 0.3244871420261927
$ws.Range("P8").Value = 0.3244871420261927
$ws.Range("S8").Value = 0.044434451217105
$ws.Range("T8").Value = 0.04443445121710499
$ws.Range("I9").Value = 0.1369374790620155
$ws.Range("J9").Value = 0.1369374790620154
$ws.Range("M9").Value = 3.197710666666667
$ws.Range("N9").Value = 9.593132000000001
$ws.Range("O9").Value = 0.08329089836363292
$ws.Range("P9").Value = 0.0832908983636329
$ws.Range("Q9").Value = 255.2644573428951
$ws.Range("R9").Value = 2297.380116086056
$ws.Range("S9").Value = 0.01140564565072644
$ws.Range("T9").Value = 0.01140564565072644
$ws.Range("I10").Value = 0.1369374790620155
$ws.Range("J10").Value = 0.1369374790620154
$ws.Range("M10").Value = 8.081220666666667
$ws.Range("N10").Value = 24.243662
$ws.Range("O10").Value = 0.2104918797744333
$ws.Range("P10").Value = 0.2104918797744333
$ws.Range("Q10").Value = 645.1016440130885
$ws.Range("R10").Value = 5805.914796117796
$ws.Range("S10").Value = 0.02882422737933574
$ws.Range("T10").Value = 0.02882422737933574
$ws.Range("I11").Value = 0.1369374790620155
$ws.Range("J11").Value = 0.1369374790620154
$ws.Range("M11").Value = 2.844675333333333
$ws.Range("N11").Value = 8.534026000000001
$ws.Range("O11").Value = 0.07409537283533685
$ws.Range("P11").Value = 0.07409537283533686
$ws.Range("Q11").Value = 227.0826165886342
$ws.Range("R11").Value = 2043.743549297708
$ws.Range("S11").Value = 0.01014643356623117
$ws.Range("T11").Value = 0.01014643356623117
$ws.Range("G12").Value = 148.824417
$ws.Range("H12").Value = 446.473251
$ws.Range("I12").Value = 0.2552967790580629
$ws.Range("J12").Value = 0.2552967790580629
$ws.Range("M12").Value = 11.81073566666667
$ws.Range("N12").Value = 35.432207
$ws.Range("O12").Value = 0.3076347070004043
$ws.Range("P12").Value = 0.3076347070004043
$ws.Range("Q12").Value = 1757.725849932773
$ws.Range("R12").Value = 15819.53264939496
$ws.Range("S12").Value = 0.07853814982367412
$ws.Range("T12").Value = 0.07853814982367412
$ws.Range("G13").Value = 148.824417
$ws.Range("H13").Value = 446.473251
$ws.Range("I13").Value = 0.2552967790580629
$ws.Range("J13").Value = 0.2552967790580629
$ws.Range("O13").Value = 0.3244871420261927
$ws.Range("P13").Value = 0.3244871420261927
$ws.Range("Q13").Value = 1854.015247731773
$ws.Range("R13").Value = 16686.13722958596
$ws.Range("S13").Value = 0.08284052220504318
$ws.Range("T13").Value = 0.08284052220504318
$ws.Range("G14").Value = 148.824417
$ws.Range("H14").Value = 446.473251
$ws.Range("I14").Value = 0.2552967790580629
$ws.Range("J14").Value = 0.2552967790580629
$ws.Range("M14").Value = 3.197710666666667
$ws.Range("N14").Value = 9.593132000000001
$ws.Range("O14").Value = 0.08329089836363292
$ws.Range("P14").Value = 0.0832908983636329
$ws.Range("Q14").Value = 475.8974257013481
$ws.Range("R14").Value = 4283.076831312133
$ws.Range("S14").Value = 0.02126389807708796
$ws.Range("T14").Value = 0.02126389807708796
$ws.Range("G15").Value = 148.824417
$ws.Range("H15").Value = 446.473251
$ws.Range("I15").Value = 0.2552967790580629
$ws.Range("J15").Value = 0.2552967790580629
$ws.Range("M15").Value = 8.081220666666667
$ws.Range("N15").Value = 24.243662
$ws.Range("O15").Value = 0.2104918797744333
$ws.Range("P15").Value = 0.2104918797744333
$ws.Range("Q15").Value = 1202.682954365018
$ws.Range("R15").Value = 10824.14658928516
$ws.Range("S15").Value = 0.05373789892428984
$ws.Range("T15").Value = 0.05373789892428984
$ws.Range("G16").Value = 148.824417
$ws.Range("H16").Value = 446.473251
$ws.Range("I16").Value = 0.2552967790580629
$ws.Range("J16").Value = 0.2552967790580629
$ws.Range("M16").Value = 2.844675333333333
$ws.Range("N16").Value = 8.534026000000001
$ws.Range("O16").Value = 0.07409537283533685
$ws.Range("P16").Value = 0.07409537283533686
$ws.Range("Q16").Value = 423.357148037614
$ws.Range("R16").Value = 3810.214332338526
$ws.Range("S16").Value = 0.01891631002796778
$ws.Range("T16").Value = 0.01891631002796779
$ws.Range("G17").Value = 35.426853
$ws.Range("H17").Value = 106.280559
$ws.Range("I17").Value = 0.06077202683121193
$ws.Range("J17").Value = 0.06077202683121192
$ws.Range("M17").Value = 11.81073566666667
$ws.Range("N17").Value = 35.432207
$ws.Range("O17").Value = 0.3076347070004043
$ws.Range("P17").Value = 0.3076347070004043
$ws.Range("Q17").Value = 418.417196284857
$ws.Range("R17").Value = 3765.754766563713
$ws.Range("S17").Value = 0.01869558466804059
$ws.Range("T17").Value = 0.01869558466804059
$ws.Range("G18").Value = 35.426853
$ws.Range("H18").Value = 106.280559
$ws.Range("I18").Value = 0.06077202683121193
$ws.Range("J18").Value = 0.06077202683121192
$ws.Range("O18").Value = 0.3244871420261927
$ws.Range("P18").Value = 0.3244871420261927
$ws.Range("Q18").Value = 441.338370175857
$ws.Range("R18").Value = 3972.045331582713
$ws.Range("S18").Value = 0.01971974130159906
$ws.Range("T18").Value = 0.01971974130159905
$ws.Range("G19").Value = 35.426853
$ws.Range("H19").Value = 106.280559
$ws.Range("I19").Value = 0.06077202683121193
$ws.Range("J19").Value = 0.06077202683121192
$ws.Range("M19").Value = 3.197710666666667
$ws.Range("N19").Value = 9.593132000000001
$ws.Range("O19").Value = 0.08329089836363292
$ws.Range("P19").Value = 0.0832908983636329
$ws.Range("Q19").Value = 113.284825724532
$ws.Range("R19").Value = 1019.563431520788
$ws.Range("S19").Value = 0.005061756710150445
$ws.Range("T19").Value = 0.005061756710150443
$ws.Range("G20").Value = 35.426853
$ws.Range("H20").Value = 106.280559
$ws.Range("I20").Value = 0.06077202683121193
$ws.Range("J20").Value = 0.06077202683121192
$ws.Range("M20").Value = 8.081220666666667
$ws.Range("N20").Value = 24.243662
$ws.Range("O20").Value = 0.2104918797744333
$ws.Range("P20").Value = 0.2104918797744333
$ws.Range("Q20").Value = 286.292216618562
$ws.Range("R20").Value = 2576.629949567058
$ws.Range("S20").Value = 0.0127920181654041
$ws.Range("T20").Value = 0.0127920181654041
$ws.Range("G21").Value = 35.426853
$ws.Range("H21").Value = 106.280559
$ws.Range("I21").Value = 0.06077202683121193
$ws.Range("J21").Value = 0.06077202683121192
$ws.Range("M21").Value = 2.844675333333333
$ws.Range("N21").Value = 8.534026000000001
$ws.Range("O21").Value = 0.07409537283533685
$ws.Range("P21").Value = 0.07409537283533686
$ws.Range("Q21").Value = 100.777894866726
$ws.Range("R21").Value = 907.001053800534
$ws.Range("S21").Value = 0.004502925986017742
$ws.Range("T21").Value = 0.004502925986017743
$ws.Range("G22").Value = 121.3248153333333
$ws.Range("H22").Value = 363.974446
$ws.Range("I22").Value = 0.2081233388901116
$ws.Range("J22").Value = 0.2081233388901115
$ws.Range("M22").Value = 11.81073566666667
$ws.Range("N22").Value = 35.432207
$ws.Range("O22").Value = 0.3076347070004043
$ws.Range("P22").Value = 0.3076347070004043
$ws.Range("Q22").Value = 1432.935323709147
$ws.Range("R22").Value = 12896.41791338232
$ws.Range("S22").Value = 0.06402596237940532
$ws.Range("T22").Value = 0.0640259623794053
$ws.Range("G23").Value = 121.3248153333333
$ws.Range("H23").Value = 363.974446
$ws.Range("I23").Value = 0.2081233388901116
$ws.Range("J23").Value = 0.2081233388901115
$ws.Range("O23").Value = 0.3244871420261927
$ws.Range("P23").Value = 0.3244871420261927
$ws.Range("Q23").Value = 1511.432479229813
$ws.Range("R23").Value = 13602.89231306832
$ws.Range("S23").Value = 0.06753334742540106
$ws.Range("T23").Value = 0.06753334742540104
$ws.Range("G24").Value = 121.3248153333333
$ws.Range("H24").Value = 363.974446
$ws.Range("I24").Value = 0.2081233388901116
$ws.Range("J24").Value = 0.2081233388901115
$ws.Range("M24").Value = 3.197710666666667
$ws.Range("N24").Value = 9.593132000000001
$ws.Range("O24").Value = 0.08329089836363292
$ws.Range("P24").Value = 0.0832908983636329
$ws.Range("Q24").Value = 387.9616561227636
$ws.Range("R24").Value = 3491.654905104872
$ws.Range("S24").Value = 0.01733477986659621
$ws.Range("T24").Value = 0.01733477986659621
$ws.Range("G25").Value = 121.3248153333333
$ws.Range("H25").Value = 363.974446
$ws.Range("I25").Value = 0.2081233388901116
$ws.Range("J25").Value = 0.2081233388901115
$ws.Range("M25").Value = 8.081220666666667
$ws.Range("N25").Value = 24.243662
$ws.Range("O25").Value = 0.2104918797744333
$ws.Range("P25").Value = 0.2104918797744333
$ws.Range("Q25").Value = 980.4526050512502
$ws.Range("R25").Value = 8824.073445461252
$ws.Range("S25").Value = 0.04380827282791101
$ws.Range("T25").Value = 0.04380827282791101
$ws.Range("G26").Value = 121.3248153333333
$ws.Range("H26").Value = 363.974446
$ws.Range("I26").Value = 0.2081233388901116
$ws.Range("J26").Value = 0.2081233388901115
$ws.Range("M26").Value = 2.844675333333333
$ws.Range("N26").Value = 8.534026000000001
$ws.Range("O26").Value = 0.07409537283533685
$ws.Range("P26").Value = 0.07409537283533686
$ws.Range("Q26").Value = 345.1297094999551
$ws.Range("R26").Value = 3106.167385499596
$ws.Range("S26").Value = 0.01542097639079798
$ws.Range("T26").Value = 0.01542097639079798
